$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (02-January-2025): move the "At Work" day to "Sick Leave"
$ws.Range("C12").Value = $null
$ws.Range("E12").Value = 1

# Row 13 (03-January-2025): move the "At Work" day to "Sick Leave"
$ws.Range("C13").Value = $null
$ws.Range("E13").Value = 1

# Row 44 (Total): update totals to reflect the two reclassified days
$ws.Range("C44").Value = 18
$ws.Range("E44").Value = 2
